$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CP003 -> CP003_Eminent (test-case id) and its description text updated
$ws.Range("B4").Value = "Sucursales Exclusivas Éminent:"
$ws.Range("A4").Value = "CP003_Eminent"

# Move the active selection from C6 to A4
$ws.Range("A4").Select()
